# Auto-generated edit script for CompStat weekly report update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/21/2025  Through  7/27/2025"

# --- Crime Complaints table value updates ---
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 4
$ws.Range("F14").Value = 5
$ws.Range("H14").Value = -50
$ws.Range("I14").Value = 57
$ws.Range("J14").Value = 67
$ws.Range("K14").Value = -14.925373134328
$ws.Range("L14").Value = -26.923076923076
$ws.Range("M14").Value = -24
$ws.Range("N14").Value = -79.787234042553

# Row 15
$ws.Range("F15").Value = 52
$ws.Range("G15").Value = 24
$ws.Range("H15").Value = 116.666666666667
$ws.Range("I15").Value = 310
$ws.Range("J15").Value = 240
$ws.Range("K15").Value = 29.166666666666
$ws.Range("L15").Value = 36.563876651982
$ws.Range("M15").Value = 87.878787878787
$ws.Range("N15").Value = -24.019607843137

# Row 16
$ws.Range("C16").Value = 101
$ws.Range("D16").Value = 106
$ws.Range("E16").Value = -4.716981132075
$ws.Range("F16").Value = 437
$ws.Range("G16").Value = 439
$ws.Range("H16").Value = -0.455580865603
$ws.Range("I16").Value = 2722
$ws.Range("J16").Value = 2818
$ws.Range("K16").Value = -3.406671398154
$ws.Range("L16").Value = 1.605076521089
$ws.Range("M16").Value = 11.648892534864
$ws.Range("N16").Value = -69.975733509816

# Row 17
$ws.Range("C17").Value = 181
$ws.Range("D17").Value = 194
$ws.Range("E17").Value = -6.701030927835
$ws.Range("F17").Value = 740
$ws.Range("G17").Value = 751
$ws.Range("H17").Value = -1.464713715046
$ws.Range("I17").Value = 4992
$ws.Range("J17").Value = 4760
$ws.Range("K17").Value = 4.873949579831
$ws.Range("L17").Value = 8.710801393728
$ws.Range("M17").Value = 98.173878523223
$ws.Range("N17").Value = -4.147465437788

# Row 18
$ws.Range("C18").Value = 63
$ws.Range("D18").Value = 74
$ws.Range("E18").Value = -14.864864864864
$ws.Range("F18").Value = 210
$ws.Range("G18").Value = 250
$ws.Range("H18").Value = -16
$ws.Range("I18").Value = 1638
$ws.Range("J18").Value = 1674
$ws.Range("K18").Value = -2.150537634408
$ws.Range("L18").Value = -4.767441860465
$ws.Range("M18").Value = -9.502762430939
$ws.Range("N18").Value = -84.655737704918

# Row 19
$ws.Range("C19").Value = 170
$ws.Range("D19").Value = 210
$ws.Range("E19").Value = -19.047619047619
$ws.Range("F19").Value = 743
$ws.Range("G19").Value = 771
$ws.Range("H19").Value = -3.631647211413
$ws.Range("I19").Value = 5160
$ws.Range("J19").Value = 5203
$ws.Range("K19").Value = -0.826446280991
$ws.Range("L19").Value = 16.425992779783
$ws.Range("M19").Value = 101.720093823299
$ws.Range("N19").Value = 24.697921701305

# Row 20
$ws.Range("C20").Value = 97
$ws.Range("D20").Value = 81
$ws.Range("E20").Value = 19.753086419753
$ws.Range("F20").Value = 385
$ws.Range("H20").Value = 5.769230769230
$ws.Range("I20").Value = 2566
$ws.Range("J20").Value = 2348
$ws.Range("K20").Value = 9.284497444633
$ws.Range("L20").Value = -15.647600262984
$ws.Range("M20").Value = 121.588946459413
$ws.Range("N20").Value = -70.441193410897

# Row 21
$ws.Range("C21").Value = 627
$ws.Range("D21").Value = 675
$ws.Range("E21").Value = -7.111111111111
$ws.Range("F21").Value = 2572
$ws.Range("G21").Value = 2609
$ws.Range("H21").Value = -1.418167880413
$ws.Range("I21").Value = 17445
$ws.Range("J21").Value = 17110
$ws.Range("K21").Value = 1.957919345412
$ws.Range("L21").Value = 4.025044722719
$ws.Range("M21").Value = 62.687680686375
$ws.Range("N21").Value = -54.638826772063

# Row 22
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -57.142857142857
$ws.Range("I22").Value = 171
$ws.Range("J22").Value = 186
$ws.Range("K22").Value = -8.064516129032
$ws.Range("L22").Value = 2.395209580838
$ws.Range("M22").Value = -7.567567567567

# Row 23
$ws.Range("C23").Value = 34
$ws.Range("D23").Value = 40
$ws.Range("E23").Value = -15
$ws.Range("F23").Value = 116
$ws.Range("G23").Value = 150
$ws.Range("H23").Value = -22.666666666666
$ws.Range("I23").Value = 893
$ws.Range("J23").Value = 982
$ws.Range("K23").Value = -9.063136456211
$ws.Range("L23").Value = -12.536728697355
$ws.Range("M23").Value = 47.847682119205

# Row 24
$ws.Range("C24").Value = 475
$ws.Range("D24").Value = 367
$ws.Range("E24").Value = 29.427792915531
$ws.Range("F24").Value = 1564
$ws.Range("G24").Value = 1250
$ws.Range("H24").Value = 25.12
$ws.Range("I24").Value = 10297
$ws.Range("J24").Value = 9141
$ws.Range("K24").Value = 12.646318783502
$ws.Range("L24").Value = 0.595935912465
$ws.Range("M24").Value = 44.661421747681

# Row 25
$ws.Range("C25").Value = 169
$ws.Range("D25").Value = 144
$ws.Range("E25").Value = 17.361111111111
$ws.Range("F25").Value = 559
$ws.Range("G25").Value = 490
$ws.Range("H25").Value = 14.081632653061
$ws.Range("I25").Value = 3426
$ws.Range("J25").Value = 3650
$ws.Range("K25").Value = -6.136986301369
$ws.Range("L25").Value = -22.242396731729

# Row 26
$ws.Range("C26").Value = 233
$ws.Range("D26").Value = 233
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 980
$ws.Range("G26").Value = 894
$ws.Range("H26").Value = 9.619686800894
$ws.Range("I26").Value = 6394
$ws.Range("J26").Value = 6366
$ws.Range("K26").Value = 0.439836632108
$ws.Range("L26").Value = 6.762397729170
$ws.Range("M26").Value = 0.172332758890

# Row 27
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = 12
$ws.Range("E27").Value = 16.666666666666
$ws.Range("F27").Value = 59
$ws.Range("G27").Value = 41
$ws.Range("H27").Value = 43.902439024390
$ws.Range("I27").Value = 387
$ws.Range("J27").Value = 378
$ws.Range("K27").Value = 2.380952380952
$ws.Range("L27").Value = 3.2

# Row 28
$ws.Range("C28").Value = 20
$ws.Range("D28").Value = 29
$ws.Range("E28").Value = -31.034482758620
$ws.Range("F28").Value = 83
$ws.Range("G28").Value = 89
$ws.Range("H28").Value = -6.741573033707
$ws.Range("I28").Value = 627
$ws.Range("J28").Value = 696
$ws.Range("K28").Value = -9.913793103448
$ws.Range("L28").Value = 4.326123128119

# Row 29
$ws.Range("C29").Value = 12
$ws.Range("D29").Value = 17
$ws.Range("E29").Value = -29.411764705882
$ws.Range("F29").Value = 36
$ws.Range("G29").Value = 48
$ws.Range("H29").Value = -25
$ws.Range("I29").Value = 181
$ws.Range("J29").Value = 233
$ws.Range("K29").Value = -22.317596566523
$ws.Range("L29").Value = -22.649572649572
$ws.Range("M29").Value = -31.439393939393
$ws.Range("N29").Value = -76.794871794871

# Row 30
$ws.Range("C30").Value = 10
$ws.Range("D30").Value = 12
$ws.Range("E30").Value = -16.666666666666
$ws.Range("G30").Value = 37
$ws.Range("H30").Value = -18.918918918918
$ws.Range("I30").Value = 157
$ws.Range("J30").Value = 188
$ws.Range("K30").Value = -16.489361702127
$ws.Range("L30").Value = -16.931216931216
$ws.Range("M30").Value = -28.959276018099
$ws.Range("N30").Value = -77.887323943662

# Row 31
$ws.Range("F31").Value = 2
$ws.Range("I31").Value = 14
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 27.272727272727

# Row 33
$ws.Range("D33").Value = 3
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = -83.333333333333
$ws.Range("J33").Value = 29
$ws.Range("K33").Value = -51.724137931034
